$d = $word.ActiveDocument

# --- Step 1: the stray empty paragraph that holds only the "last edit"
#     _GoBack bookmark gets removed (its paragraph mark merges into the
#     previous, now-final, empty paragraph). ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $goBackPara = $d.Bookmarks.Item("_GoBack").Range.Paragraphs.Item(1)
    $goBackPara.Range.Delete()
}

# --- Step 2: the trailing "Lilly Customer Meeting Service Team ..." run
#     grows from 8pt to 9pt (sz/szCs 16 -> 18), and the _GoBack bookmark
#     is re-created right after that run (Word's "last edit location"
#     marker follows the edit). ---
$findRange = $d.Content
$findRange.Find.Execute("Lilly Customer Meeting Service Team", `
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$runEnd = $lastPara.Range.End - 1   # paragraph end, minus the trailing pilcrow
$targetRange = $d.Range($findRange.Start, $runEnd)

$snippet = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:body><w:p>' +
  '<w:r w:rsidR="005128A9" w:rsidRPr="00D650C6">' +
  '<w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:sz w:val="18"/><w:szCs w:val="18"/><w:lang w:val="de-DE"/></w:rPr>' +
  '<w:t>Lilly Customer Meeting Service Team c/o Europe Convention GmbH &amp; Co. KG Bahnhofstra&#223;e 30 82467 Garmisch-Partenkirchen</w:t>' +
  '</w:r>' +
  '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
  '</w:p></w:body></w:document>' +
  '</pkg:xmlData></pkg:part></pkg:package>'

$targetRange.InsertXML($snippet)
Write-Host "Edit applied."
